# Applies Scott's MIP changes: adds an index/legend column on Sheet 2,
# recodes several COVID-19 rows to "Don't know" (18), removes the stray
# blank "other_problem_text" row, and adds variable labels (category
# code + label) to the previously-unlabelled dichotomous-variable rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 2")

# ---------------------------------------------------------------------
# 1. New column I (rows 1-19): index/legend of category codes.
# ---------------------------------------------------------------------
# Shared-string insertion order matters (new strings are appended to the
# table in first-use order), so write the numbered legend entries (I2:I19)
# before the "index" header (I1) to match the upstream string ordering.
$indexLabels = @(
  "1 = COVID-19",
  "2 = Public health messaging/gov't handling of COVID",
  "3 = Healthcare (access to care, short supply)",
  "4 = Long term care",
  "5 = Chronic disease (cancer, heart disease)",
  "6 = Mental health",
  "7 = Access to housing and food",
  "8 = Drug abuse",
  "9 = Inequality",
  "10 = Economy",
  "11 = Corrupt gov't",
  "12 = Climate change/environmental",
  "13 = Abortion",
  "14 = Reliance on meat",
  "15 = Domestic abuse",
  "16 = Misinformation",
  "17 = Internet addiction",
  "18 = Don't know"
)

for ($i = 0; $i -lt $indexLabels.Length; $i++) {
    $ws.Cells.Item($i + 2, 9).Value = $indexLabels[$i]
}
$ws.Cells.Item(1, 9).Value = "index"

# ---------------------------------------------------------------------
# 2. Recode several "COVID-19" (1) rows to "Don't know" (18).
# ---------------------------------------------------------------------
$ws.Cells.Item(78, 2).Value = 18
$ws.Cells.Item(78, 3).Value = 18
$ws.Cells.Item(83, 3).Value = 18
$ws.Cells.Item(84, 3).Value = 18
$ws.Cells.Item(87, 3).Value = 18
$ws.Cells.Item(95, 3).Value = 18

# ---------------------------------------------------------------------
# 3. Remove the stray blank "other_problem_text" row (old row 107),
#    shifting the remaining rows up.
# ---------------------------------------------------------------------
$ws.Rows.Item(107).Delete()

# ---------------------------------------------------------------------
# 4. Add variable labels (category code + label) to the previously
#    unlabelled dichotomous-variable rows (now rows 107-114).
# ---------------------------------------------------------------------
$ws.Cells.Item(107, 3).Value = 6
$ws.Cells.Item(107, 4).Value = "Mental health"

$ws.Cells.Item(108, 3).Value = 2
$ws.Cells.Item(108, 4).Value = "Public health messaging / gov't handing of COVID"

$ws.Cells.Item(109, 3).Value = 5
$ws.Cells.Item(109, 4).Value = "Chronic disease (cancer, heart disease, etc)"

$ws.Cells.Item(110, 3).Value = 7
$ws.Cells.Item(110, 4).Value = "Access to housing and food"

$ws.Cells.Item(111, 3).Value = 7
$ws.Cells.Item(111, 4).Value = "Access to housing and food"

$ws.Cells.Item(112, 3).Value = 7
$ws.Cells.Item(112, 4).Value = "Access to housing and food"

$ws.Cells.Item(113, 3).Value = 2
$ws.Cells.Item(113, 4).Value = "Mental health"

$ws.Cells.Item(114, 3).Value = 9
$ws.Cells.Item(114, 4).Value = "Inequality and discrimination"

$wb.Save()
